# Update CBS model excess mortality table
# - Inserts a new "Jaar" (Year) column before the existing week-number column.
# - Adds year values to existing data rows.
# - Corrects several "Waargenomen" (observed) values.
# - Appends week 53/2020 and week 1/2021 data rows.
# - Moves / rebuilds the totals row further down with updated sums.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new column before column F. Excel shifts F->G, G->H, H->I,
#    I->J automatically (including formula references).
# ---------------------------------------------------------------------
$ws.Columns("F").Insert()

# Match the new column's width to its neighbour (column E).
$ws.Columns("F").ColumnWidth = $ws.Columns("E").ColumnWidth

# ---------------------------------------------------------------------
# 2. Header row: F1 = "Jaar", G1 = "Week" (existing G1.."Waargenomen" etc.
#    already shifted right to H1/I1/J1 by the column insert above).
#    ("Week" is written first so it lands before "Jaar" in the shared
#    string table, matching the source workbook's string order.)
# ---------------------------------------------------------------------
$ws.Range("G1").Value = "Week"
$ws.Range("F1").Value = "Jaar"

# ---------------------------------------------------------------------
# 3. Per-row data: year (F), week (G, already correct post-shift),
#    corrected "Waargenomen" (H) and "Verwacht" (I) values.
#    G / the week number does not change - only listed for clarity/verification.
# ---------------------------------------------------------------------
$rows = @(
    @{Row=3; Year=2020; Week=11; H=3219; I=3253},
    @{Row=4; Year=2020; Week=12; H=3615; I=3174},
    @{Row=5; Year=2020; Week=13; H=4459; I=3104},
    @{Row=6; Year=2020; Week=14; H=5085; I=3024},
    @{Row=7; Year=2020; Week=15; H=4981; I=2957},
    @{Row=8; Year=2020; Week=16; H=4307; I=2915},
    @{Row=9; Year=2020; Week=17; H=3907; I=2869},
    @{Row=10; Year=2020; Week=18; H=3379; I=2841},
    @{Row=11; Year=2020; Week=19; H=2986; I=2821},
    @{Row=12; Year=2020; Week=20; H=2777; I=2794},
    @{Row=13; Year=2020; Week=21; H=2771; I=2770},
    @{Row=14; Year=2020; Week=22; H=2729; I=2753},
    @{Row=15; Year=2020; Week=23; H=2682; I=2735},
    @{Row=16; Year=2020; Week=24; H=2692; I=2737},
    @{Row=17; Year=2020; Week=25; H=2695; I=2725},
    @{Row=18; Year=2020; Week=26; H=2661; I=2717},
    @{Row=19; Year=2020; Week=27; H=2639; I=2723},
    @{Row=20; Year=2020; Week=28; H=2619; I=2719},
    @{Row=21; Year=2020; Week=29; H=2528; I=2720},
    @{Row=22; Year=2020; Week=30; H=2673; I=2707},
    @{Row=23; Year=2020; Week=31; H=2668; I=2687},
    @{Row=24; Year=2020; Week=32; H=2640; I=2682},
    @{Row=25; Year=2020; Week=33; H=3209; I=2669},
    @{Row=26; Year=2020; Week=34; H=2855; I=2663},
    @{Row=27; Year=2020; Week=35; H=2733; I=2667},
    @{Row=28; Year=2020; Week=36; H=2690; I=2676},
    @{Row=29; Year=2020; Week=37; H=2739; I=2698},
    @{Row=30; Year=2020; Week=38; H=2720; I=2729},
    @{Row=31; Year=2020; Week=39; H=2892; I=2752},
    @{Row=32; Year=2020; Week=40; H=2998; I=2786},
    @{Row=33; Year=2020; Week=41; H=3020; I=2807},
    @{Row=34; Year=2020; Week=42; H=3220; I=2839},
    @{Row=35; Year=2020; Week=43; H=3449; I=2862},
    @{Row=36; Year=2020; Week=44; H=3679; I=2889},
    @{Row=37; Year=2020; Week=45; H=3589; I=2902},
    @{Row=38; Year=2020; Week=46; H=3575; I=2932},
    @{Row=39; Year=2020; Week=47; H=3329; I=2972},
    @{Row=40; Year=2020; Week=48; H=3401; I=3012},
    @{Row=41; Year=2020; Week=49; H=3519; I=3037},
    @{Row=42; Year=2020; Week=50; H=3606; I=3100},
    @{Row=43; Year=2020; Week=51; H=3896; I=3166},
    @{Row=44; Year=2020; Week=52; H=3849; I=3222},
    @{Row=45; Year=2020; Week=53; H=4058; I=3266},
    @{Row=46; Year=2021; Week=1;  H=3954; I=3309}
)

# Row 2 (the "3-10" average row) also gets a year value.
$ws.Range("F2").Value = 2020

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 6).Value = $r.Year     # F - Jaar
    $ws.Cells.Item($row, 7).Value = $r.Week     # G - Week
    $ws.Cells.Item($row, 8).Value = $r.H        # H - Waargenomen
    $ws.Cells.Item($row, 9).Value = $r.I        # I - Verwacht
    $ws.Cells.Item($row, 10).Formula = "=H$row-I$row"   # J - Oversterfte
}

# ---------------------------------------------------------------------
# 4. Totals row moves from row 46 to row 50, with the updated sum ranges.
#    (The loop above already overwrote the old label/formulas that used
#    to live on row 46 with real week-53/2020 + week-1/2021 data.)
# ---------------------------------------------------------------------
$ws.Range("G50").Value = "Som week 11 tot en met 19"
$ws.Range("H50").Formula = "=SUM(H3:H28)"
$ws.Range("I50").Formula = "=SUM(I3:I28)"
$ws.Range("J50").Formula = "=SUM(J3:J34)"

# ---------------------------------------------------------------------
# 5. View / selection cosmetics from the diff (scroll position + the
#    newly active cell/selection).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F47").Select()
